$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> updated D (Price) and/or E (Volume 1h) values,
# taken from the refreshed cryptos list for this run.
$changes = @{
    2 = @{ D = '27.447.31'; E = '  -0.99%  ' }
    3 = @{ D = '1.734.27'; E = '  -1.16%  ' }
    4 = @{ E = '  -0.65%  ' }
    5 = @{ D = '322.63'; E = '  -0.35%  ' }
    6 = @{ E = '  -0.45%  ' }
    7 = @{ D = '0.4534'; E = '  +6.57%  ' }
    8 = @{ D = '0.3530'; E = '  -2.62%  ' }
    9 = @{ D = '0.07372'; E = '  -2.52%  ' }
    10 = @{ D = '41.31'; E = '  -3.00%  ' }
    11 = @{ E = '  -1.59%  ' }
    12 = @{ D = '1.000'; E = '  -0.37%  ' }
    13 = @{ E = '  -1.03%  ' }
    14 = @{ D = '5.920'; E = '  -2.02%  ' }
    15 = @{ D = '7.065'; E = '  -2.74%  ' }
    16 = @{ D = '1.726.78'; E = '  -3.07%  ' }
    17 = @{ D = '90.96'; E = '  -0.32%  ' }
    18 = @{ E = '  -2.22%  ' }
    19 = @{ E = '  -0.89%  ' }
    20 = @{ D = '1.001'; E = '  -0.20%  ' }
    21 = @{ D = '16.61' }
    22 = @{ D = '5.732'; E = '  -2.97%  ' }
    23 = @{ D = '27.483.42'; E = '  -1.09%  ' }
    24 = @{ E = '  -1.09%  ' }
    25 = @{ D = '2.060'; E = '  -2.07%  ' }
    26 = @{ D = '161.47'; E = '  +0.62%  ' }
    27 = @{ D = '19.89'; E = '  -2.10%  ' }
    28 = @{ D = '1.925.07'; E = '  -2.64%  ' }
    29 = @{ D = '2.045'; E = '  -4.21%  ' }
    30 = @{ D = '124.41'; E = '  -0.76%  ' }
    31 = @{ D = '1.040'; E = '  -6.69%  ' }
    32 = @{ D = '0.09144'; E = '  +2.81%  ' }
    33 = @{ D = '3.645'; E = '  -0.95%  ' }
    34 = @{ D = '5.374'; E = '  -3.56%  ' }
    35 = @{ D = '0.02267'; E = '  -1.42%  ' }
    36 = @{ D = '11.58'; E = '  -5.39%  ' }
    37 = @{ D = '0.05944'; E = '  -1.26%  ' }
    38 = @{ D = '0.2045'; E = '  -2.88%  ' }
    39 = @{ D = '0.6230'; E = '  -1.68%  ' }
    40 = @{ D = '4.853'; E = '  -2.21%  ' }
    41 = @{ D = '1.192'; E = '  +0.87%  ' }
    42 = @{ D = '1.367'; E = '  -1.93%  ' }
    43 = @{ D = '7.709'; E = '  -2.40%  ' }
    44 = @{ D = '13.00'; E = '  -2.63%  ' }
    45 = @{ E = '  -0.06%  ' }
    46 = @{ D = '0.5789'; E = '  -1.29%  ' }
    47 = @{ D = '121.90'; E = '  -0.87%  ' }
    48 = @{ D = '1.920'; E = '  -3.28%  ' }
    49 = @{ E = '  -0.03%  ' }
    50 = @{ D = '1.110'; E = '  -5.62%  ' }
    51 = @{ D = '70.81'; E = '  -3.68%  ' }
}

foreach ($rowNum in $changes.Keys) {
    $vals = $changes[$rowNum]
    if ($vals.ContainsKey("D")) {
        $cell = $ws.Range("D$rowNum")
        # Force text interpretation so numeric-looking strings (e.g. "1.000",
        # "0.4534") keep their exact original formatting instead of being
        # auto-converted to a Double by Excel's smart-entry parsing (which would
        # also introduce floating point noise, e.g. 0.4534 -> 0.45340000000000003).
        $cell.NumberFormat = "@"
        $cell.Value = $vals["D"]
        # Drop the temporary text format again so the cell ends up styled
        # exactly as it started (no explicit style / General format).
        $cell.ClearFormats()
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$rowNum").Value = $vals["E"]
    }
}
